$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45205 -> 45206) for every data row (header on row 1, data rows 2..239).
$lastRow = $ws.UsedRange.Rows.Count - 1
if ($lastRow -lt 239) {
    $lastRow = 239
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
